# Applies the cryptos list refresh described by the commit diff.
# All target cells are plain text cells (inline/shared strings), and several
# hold numeric-looking text (e.g. "514.47", "57.267.65") that must stay text
# rather than be auto-converted to a number by Excel. To force that, each
# cell is temporarily switched to the Text number format before the value is
# written, then its style is restored to "Normal" so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "57.267.65"
Set-TextValue $ws "E2" "  +2.93%  "

Set-TextValue $ws "D3" "3.068.73"
Set-TextValue $ws "E3" "  +5.21%  "

Set-TextValue $ws "E4" "  +0.16%  "

Set-TextValue $ws "D5" "514.47"
Set-TextValue $ws "E5" "  +1.16%  "

Set-TextValue $ws "D6" "141.01"
Set-TextValue $ws "E6" "  +5.32%  "

Set-TextValue $ws "E7" "  +0.14%  "

Set-TextValue $ws "D8" "0.434"
Set-TextValue $ws "E8" "  +3.31%  "

Set-TextValue $ws "D9" "7.24"
Set-TextValue $ws "E9" "  +0.84%  "

Set-TextValue $ws "E10" "  +3.22%  "

Set-TextValue $ws "E11" "  +5.48%  "

Set-TextValue $ws "D12" "3.595.50"
Set-TextValue $ws "E12" "  +5.52%  "

Set-TextValue $ws "E13" "  +2.63%  "

Set-TextValue $ws "D14" "25.46"
Set-TextValue $ws "E14" "  -1.64%  "

Set-TextValue $ws "E15" "  +3.06%  "

Set-TextValue $ws "D16" "57.312.91"
Set-TextValue $ws "E16" "  +3.22%  "

Set-TextValue $ws "D17" "3.069.64"
Set-TextValue $ws "E17" "  +5.43%  "

Set-TextValue $ws "D18" "5.91"
Set-TextValue $ws "E18" "  -2.22%  "

Set-TextValue $ws "D19" "13.03"
Set-TextValue $ws "E19" "  +3.44%  "

Set-TextValue $ws "D20" "8.13"
Set-TextValue $ws "E20" "  +5.95%  "

Set-TextValue $ws "D21" "336.76"
Set-TextValue $ws "E21" "  +6.59%  "

Set-TextValue $ws "E22" "  +0.07%  "

Set-TextValue $ws "D23" "0.500"
Set-TextValue $ws "E23" "  +3.00%  "

Set-TextValue $ws "D24" "65.37"
Set-TextValue $ws "E24" "  +4.50%  "

Set-TextValue $ws "D25" "0.170"
Set-TextValue $ws "E25" "  +6.22%  "

Set-TextValue $ws "E26" "  +0.31%  "

Set-TextValue $ws "D27" "0.0₃0947"
Set-TextValue $ws "E27" "  +11.74%  "

Set-TextValue $ws "D28" "6.44"
Set-TextValue $ws "E28" "  +0.83%  "

Set-TextValue $ws "D29" "7.05"
Set-TextValue $ws "E29" "  +1.58%  "

Set-TextValue $ws "E30" "  +2.07%  "

Set-TextValue $ws "D31" "20.72"
Set-TextValue $ws "E31" "  +4.71%  "

Set-TextValue $ws "D32" "1.18"
Set-TextValue $ws "E32" "  +3.88%  "

Set-TextValue $ws "D33" "154.02"
Set-TextValue $ws "E33" "  +3.59%  "

Set-TextValue $ws "D34" "4.54"
Set-TextValue $ws "E34" "  +3.21%  "

Set-TextValue $ws "D35" "5.86"
Set-TextValue $ws "E35" "  +4.29%  "

Set-TextValue $ws "D36" "25.98"
Set-TextValue $ws "E36" "  +5.69%  "

Set-TextValue $ws "D37" "1.24"
Set-TextValue $ws "E37" "  +4.15%  "

Set-TextValue $ws "D38" "0.0671"
Set-TextValue $ws "E38" "  +3.05%  "

Set-TextValue $ws "D39" "3.106.60"
Set-TextValue $ws "E39" "  +5.50%  "

Set-TextValue $ws "D40" "37.02"
Set-TextValue $ws "E40" "  +2.06%  "

Set-TextValue $ws "D41" "0.670"
Set-TextValue $ws "E41" "  +5.37%  "

Set-TextValue $ws "B42" "Filecoin"
Set-TextValue $ws "C42" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws "D42" "3.84"
Set-TextValue $ws "E42" "  +3.52%  "

Set-TextValue $ws "B43" "FirstDigitalUSD"
Set-TextValue $ws "C43" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws "D43" "1.00"
Set-TextValue $ws "E43" "  +0.27%  "

Set-TextValue $ws "D44" "2.241.79"
Set-TextValue $ws "E44" "  +6.69%  "

Set-TextValue $ws "E45" "  +8.45%  "

Set-TextValue $ws "D46" "1.38"
Set-TextValue $ws "E46" "  +3.72%  "

Set-TextValue $ws "D47" "0.948"
Set-TextValue $ws "E47" "  +3.68%  "

Set-TextValue $ws "D48" "19.98"
Set-TextValue $ws "E48" "  +7.13%  "

Set-TextValue $ws "E49" "  -1.38%  "

Set-TextValue $ws "D50" "0.0866"
Set-TextValue $ws "E50" "  +3.59%  "

Set-TextValue $ws "B51" "SuiNetwork"
Set-TextValue $ws "C51" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws "D51" "0.687"
Set-TextValue $ws "E51" "  +5.49%  "

